$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-10-18 01:39:39"
}
